# Regenerate the memory/bedrooms input list (block 2) after the author
# finished input-list generation + sanity checks: trial_total (F) is
# renumbered and every trial row (2-29) has its condition / stimulus /
# pretest-norm columns (G:S) reassigned to the corrected stimulus set.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# row 2 (trial_block 1)
$ws.Cells.Item(2, 6).Value = 189  # F2 trial_total
$ws.Cells.Item(2, 9).Value = 'target'  # I2 cond_cat
$ws.Cells.Item(2, 10).Value = 'old'  # J2 cond_mem
$ws.Cells.Item(2, 11).Value = 'j'  # K2 correct_answer
$ws.Cells.Item(2, 12).Value = 'stimuli/img_le8uf.png'  # L2 stimulus
$ws.Cells.Item(2, 13).Value = 12.88888888888889  # M2 conceptual
$ws.Cells.Item(2, 14).Value = 9.222222222222221  # N2 perceptual
$ws.Cells.Item(2, 15).Value = 11.05555555555556  # O2 typicality
$ws.Cells.Item(2, 16).Value = 36  # P2 n
$ws.Cells.Item(2, 17).Value = 1  # Q2 p_typicality
$ws.Cells.Item(2, 18).Value = 1  # R2 p_conceptual
$ws.Cells.Item(2, 19).Value = 1  # S2 p_perceptual

# row 3 (trial_block 2)
$ws.Cells.Item(3, 6).Value = 190  # F3 trial_total
$ws.Cells.Item(3, 9).Value = 'target'  # I3 cond_cat
$ws.Cells.Item(3, 10).Value = 'old'  # J3 cond_mem
$ws.Cells.Item(3, 11).Value = 'j'  # K3 correct_answer
$ws.Cells.Item(3, 12).Value = 'stimuli/img_v8dra.png'  # L3 stimulus
$ws.Cells.Item(3, 13).Value = 61.77272727272727  # M3 conceptual
$ws.Cells.Item(3, 14).Value = 38.79545454545455  # N3 perceptual
$ws.Cells.Item(3, 15).Value = 50.28409090909091  # O3 typicality
$ws.Cells.Item(3, 16).Value = 44  # P3 n
$ws.Cells.Item(3, 17).Value = 3  # Q3 p_typicality
$ws.Cells.Item(3, 18).Value = 3  # R3 p_conceptual
$ws.Cells.Item(3, 19).Value = 3  # S3 p_perceptual

# row 4 (trial_block 3)
$ws.Cells.Item(4, 6).Value = 191  # F4 trial_total
$ws.Cells.Item(4, 12).Value = 'stimuli/img_okvvw.png'  # L4 stimulus
$ws.Cells.Item(4, 13).Value = 50.58333333333334  # M4 conceptual
$ws.Cells.Item(4, 14).Value = 32.11111111111111  # N4 perceptual
$ws.Cells.Item(4, 15).Value = 41.34722222222223  # O4 typicality
$ws.Cells.Item(4, 16).Value = 36  # P4 n
$ws.Cells.Item(4, 17).Value = 2  # Q4 p_typicality
$ws.Cells.Item(4, 18).Value = 2  # R4 p_conceptual
$ws.Cells.Item(4, 19).Value = 2  # S4 p_perceptual

# row 5 (trial_block 4)
$ws.Cells.Item(5, 6).Value = 192  # F5 trial_total
$ws.Cells.Item(5, 9).Value = 'target'  # I5 cond_cat
$ws.Cells.Item(5, 10).Value = 'old'  # J5 cond_mem
$ws.Cells.Item(5, 11).Value = 'j'  # K5 correct_answer
$ws.Cells.Item(5, 12).Value = 'stimuli/img_2js6m.png'  # L5 stimulus
$ws.Cells.Item(5, 13).Value = 40.02777777777778  # M5 conceptual
$ws.Cells.Item(5, 14).Value = 20.88888888888889  # N5 perceptual
$ws.Cells.Item(5, 15).Value = 30.45833333333334  # O5 typicality
$ws.Cells.Item(5, 16).Value = 36  # P5 n
$ws.Cells.Item(5, 17).Value = 2  # Q5 p_typicality
$ws.Cells.Item(5, 18).Value = 2  # R5 p_conceptual
$ws.Cells.Item(5, 19).Value = 2  # S5 p_perceptual

# row 6 (trial_block 5)
$ws.Cells.Item(6, 6).Value = 193  # F6 trial_total
$ws.Cells.Item(6, 12).Value = 'stimuli/img_71mhq.png'  # L6 stimulus
$ws.Cells.Item(6, 13).Value = 69.34210526315789  # M6 conceptual
$ws.Cells.Item(6, 14).Value = 47.02631578947368  # N6 perceptual
$ws.Cells.Item(6, 15).Value = 58.18421052631579  # O6 typicality
$ws.Cells.Item(6, 16).Value = 38  # P6 n
$ws.Cells.Item(6, 17).Value = 5  # Q6 p_typicality
$ws.Cells.Item(6, 18).Value = 5  # R6 p_conceptual
$ws.Cells.Item(6, 19).Value = 5  # S6 p_perceptual

# row 7 (trial_block 6)
$ws.Cells.Item(7, 6).Value = 194  # F7 trial_total
$ws.Cells.Item(7, 12).Value = 'stimuli/img_oou46.png'  # L7 stimulus
$ws.Cells.Item(7, 13).Value = 75.70270270270271  # M7 conceptual
$ws.Cells.Item(7, 14).Value = 54.86486486486486  # N7 perceptual
$ws.Cells.Item(7, 15).Value = 65.28378378378379  # O7 typicality
$ws.Cells.Item(7, 16).Value = 37  # P7 n
$ws.Cells.Item(7, 17).Value = 6  # Q7 p_typicality
$ws.Cells.Item(7, 18).Value = 6  # R7 p_conceptual
$ws.Cells.Item(7, 19).Value = 6  # S7 p_perceptual

# row 8 (trial_block 7)
$ws.Cells.Item(8, 6).Value = 195  # F8 trial_total
$ws.Cells.Item(8, 9).Value = 'target'  # I8 cond_cat
$ws.Cells.Item(8, 10).Value = 'old'  # J8 cond_mem
$ws.Cells.Item(8, 11).Value = 'j'  # K8 correct_answer
$ws.Cells.Item(8, 12).Value = 'stimuli/img_uxxo0.png'  # L8 stimulus
$ws.Cells.Item(8, 13).Value = 71.74418604651163  # M8 conceptual
$ws.Cells.Item(8, 14).Value = 48.44186046511628  # N8 perceptual
$ws.Cells.Item(8, 15).Value = 60.09302325581395  # O8 typicality
$ws.Cells.Item(8, 16).Value = 43  # P8 n
$ws.Cells.Item(8, 17).Value = 5  # Q8 p_typicality
$ws.Cells.Item(8, 18).Value = 5  # R8 p_conceptual
$ws.Cells.Item(8, 19).Value = 5  # S8 p_perceptual

# row 9 (trial_block 8)
$ws.Cells.Item(9, 6).Value = 196  # F9 trial_total
$ws.Cells.Item(9, 12).Value = 'stimuli/img_zgg62.png'  # L9 stimulus
$ws.Cells.Item(9, 13).Value = 82.18421052631579  # M9 conceptual
$ws.Cells.Item(9, 14).Value = 63.52631578947368  # N9 perceptual
$ws.Cells.Item(9, 15).Value = 72.85526315789474  # O9 typicality
$ws.Cells.Item(9, 16).Value = 38  # P9 n
$ws.Cells.Item(9, 17).Value = 8  # Q9 p_typicality
$ws.Cells.Item(9, 18).Value = 8  # R9 p_conceptual
$ws.Cells.Item(9, 19).Value = 8  # S9 p_perceptual

# row 10 (trial_block 9)
$ws.Cells.Item(10, 6).Value = 197  # F10 trial_total
$ws.Cells.Item(10, 12).Value = 'stimuli/img_5m6x4.png'  # L10 stimulus
$ws.Cells.Item(10, 13).Value = 80.23076923076923  # M10 conceptual
$ws.Cells.Item(10, 14).Value = 58.41025641025641  # N10 perceptual
$ws.Cells.Item(10, 15).Value = 69.32051282051282  # O10 typicality
$ws.Cells.Item(10, 16).Value = 39  # P10 n
$ws.Cells.Item(10, 17).Value = 7  # Q10 p_typicality
$ws.Cells.Item(10, 18).Value = 7  # R10 p_conceptual
$ws.Cells.Item(10, 19).Value = 7  # S10 p_perceptual

# row 11 (trial_block 10)
$ws.Cells.Item(11, 6).Value = 198  # F11 trial_total
$ws.Cells.Item(11, 9).Value = $null  # I11 cond_cat
$ws.Cells.Item(11, 10).Value = 'new'  # J11 cond_mem
$ws.Cells.Item(11, 11).Value = 'f'  # K11 correct_answer
$ws.Cells.Item(11, 12).Value = 'stimuli/img_zv0dq.png'  # L11 stimulus
$ws.Cells.Item(11, 13).Value = 76.86842105263158  # M11 conceptual
$ws.Cells.Item(11, 14).Value = 52.71052631578947  # N11 perceptual
$ws.Cells.Item(11, 15).Value = 64.78947368421052  # O11 typicality
$ws.Cells.Item(11, 16).Value = 38  # P11 n
$ws.Cells.Item(11, 17).Value = 6  # Q11 p_typicality
$ws.Cells.Item(11, 18).Value = 6  # R11 p_conceptual
$ws.Cells.Item(11, 19).Value = 6  # S11 p_perceptual

# row 12 (trial_block 11)
$ws.Cells.Item(12, 6).Value = 199  # F12 trial_total
$ws.Cells.Item(12, 9).Value = $null  # I12 cond_cat
$ws.Cells.Item(12, 10).Value = 'new'  # J12 cond_mem
$ws.Cells.Item(12, 11).Value = 'f'  # K12 correct_answer
$ws.Cells.Item(12, 12).Value = 'stimuli/img_6ddrx.png'  # L12 stimulus
$ws.Cells.Item(12, 13).Value = 82.2  # M12 conceptual
$ws.Cells.Item(12, 14).Value = 63.68571428571428  # N12 perceptual
$ws.Cells.Item(12, 15).Value = 72.94285714285715  # O12 typicality
$ws.Cells.Item(12, 16).Value = 35  # P12 n
$ws.Cells.Item(12, 17).Value = 8  # Q12 p_typicality
$ws.Cells.Item(12, 18).Value = 8  # R12 p_conceptual
$ws.Cells.Item(12, 19).Value = 8  # S12 p_perceptual

# row 13 (trial_block 12)
$ws.Cells.Item(13, 6).Value = 200  # F13 trial_total
$ws.Cells.Item(13, 9).Value = $null  # I13 cond_cat
$ws.Cells.Item(13, 10).Value = 'new'  # J13 cond_mem
$ws.Cells.Item(13, 11).Value = 'f'  # K13 correct_answer
$ws.Cells.Item(13, 12).Value = 'stimuli/img_0eflx.png'  # L13 stimulus
$ws.Cells.Item(13, 13).Value = 76.05128205128206  # M13 conceptual
$ws.Cells.Item(13, 14).Value = 53.53846153846154  # N13 perceptual
$ws.Cells.Item(13, 15).Value = 64.7948717948718  # O13 typicality
$ws.Cells.Item(13, 16).Value = 39  # P13 n
$ws.Cells.Item(13, 17).Value = 6  # Q13 p_typicality
$ws.Cells.Item(13, 18).Value = 6  # R13 p_conceptual
$ws.Cells.Item(13, 19).Value = 6  # S13 p_perceptual

# row 14 (trial_block 13)
$ws.Cells.Item(14, 6).Value = 201  # F14 trial_total
$ws.Cells.Item(14, 8).Value = 'bedrooms'  # H14 category
$ws.Cells.Item(14, 9).Value = 'target'  # I14 cond_cat
$ws.Cells.Item(14, 10).Value = 'old'  # J14 cond_mem
$ws.Cells.Item(14, 11).Value = 'j'  # K14 correct_answer
$ws.Cells.Item(14, 12).Value = 'stimuli/img_qgbyn.png'  # L14 stimulus
$ws.Cells.Item(14, 13).Value = 65.08108108108108  # M14 conceptual
$ws.Cells.Item(14, 14).Value = 40.10810810810811  # N14 perceptual
$ws.Cells.Item(14, 15).Value = 52.5945945945946  # O14 typicality
$ws.Cells.Item(14, 16).Value = 37  # P14 n
$ws.Cells.Item(14, 17).Value = 4  # Q14 p_typicality
$ws.Cells.Item(14, 18).Value = 4  # R14 p_conceptual
$ws.Cells.Item(14, 19).Value = 4  # S14 p_perceptual

# row 15 (trial_block 14)
$ws.Cells.Item(15, 6).Value = 202  # F15 trial_total
$ws.Cells.Item(15, 12).Value = 'stimuli/img_fqgem.png'  # L15 stimulus
$ws.Cells.Item(15, 13).Value = 80.75  # M15 conceptual
$ws.Cells.Item(15, 14).Value = 61.475  # N15 perceptual
$ws.Cells.Item(15, 15).Value = 71.1125  # O15 typicality
$ws.Cells.Item(15, 16).Value = 40  # P15 n
$ws.Cells.Item(15, 17).Value = 8  # Q15 p_typicality
$ws.Cells.Item(15, 18).Value = 8  # R15 p_conceptual
$ws.Cells.Item(15, 19).Value = 8  # S15 p_perceptual

# row 16 (trial_block 15)
$ws.Cells.Item(16, 6).Value = 203  # F16 trial_total
$ws.Cells.Item(16, 12).Value = 'stimuli/img_rvssl.png'  # L16 stimulus
$ws.Cells.Item(16, 13).Value = 74.25  # M16 conceptual
$ws.Cells.Item(16, 14).Value = 54.33333333333334  # N16 perceptual
$ws.Cells.Item(16, 15).Value = 64.29166666666667  # O16 typicality
$ws.Cells.Item(16, 16).Value = 36  # P16 n
$ws.Cells.Item(16, 17).Value = 6  # Q16 p_typicality
$ws.Cells.Item(16, 18).Value = 6  # R16 p_conceptual
$ws.Cells.Item(16, 19).Value = 6  # S16 p_perceptual

# row 17 (trial_block 16)
$ws.Cells.Item(17, 6).Value = 204  # F17 trial_total
$ws.Cells.Item(17, 12).Value = 'stimuli/img_wyctg.png'  # L17 stimulus
$ws.Cells.Item(17, 13).Value = 33.44736842105263  # M17 conceptual
$ws.Cells.Item(17, 14).Value = 11.39473684210526  # N17 perceptual
$ws.Cells.Item(17, 15).Value = 22.42105263157895  # O17 typicality
$ws.Cells.Item(17, 16).Value = 38  # P17 n
$ws.Cells.Item(17, 17).Value = 1  # Q17 p_typicality
$ws.Cells.Item(17, 18).Value = 1  # R17 p_conceptual
$ws.Cells.Item(17, 19).Value = 1  # S17 p_perceptual

# row 18 (trial_block 17)
$ws.Cells.Item(18, 6).Value = 205  # F18 trial_total
$ws.Cells.Item(18, 8).Value = $null  # H18 category
$ws.Cells.Item(18, 9).Value = $null  # I18 cond_cat
$ws.Cells.Item(18, 10).Value = 'catch'  # J18 cond_mem
$ws.Cells.Item(18, 11).Value = 'f'  # K18 correct_answer
$ws.Cells.Item(18, 12).Value = 'stimuli/catch_28.jpg'  # L18 stimulus
$ws.Cells.Item(18, 13).Value = $null  # M18 conceptual
$ws.Cells.Item(18, 14).Value = $null  # N18 perceptual
$ws.Cells.Item(18, 15).Value = $null  # O18 typicality
$ws.Cells.Item(18, 16).Value = $null  # P18 n
$ws.Cells.Item(18, 17).Value = $null  # Q18 p_typicality
$ws.Cells.Item(18, 18).Value = $null  # R18 p_conceptual
$ws.Cells.Item(18, 19).Value = $null  # S18 p_perceptual

# row 19 (trial_block 18)
$ws.Cells.Item(19, 6).Value = 206  # F19 trial_total
$ws.Cells.Item(19, 9).Value = 'target'  # I19 cond_cat
$ws.Cells.Item(19, 10).Value = 'old'  # J19 cond_mem
$ws.Cells.Item(19, 11).Value = 'j'  # K19 correct_answer
$ws.Cells.Item(19, 12).Value = 'stimuli/img_a9acb.png'  # L19 stimulus
$ws.Cells.Item(19, 13).Value = 77.11428571428571  # M19 conceptual
$ws.Cells.Item(19, 14).Value = 58.42857142857143  # N19 perceptual
$ws.Cells.Item(19, 15).Value = 67.77142857142857  # O19 typicality
$ws.Cells.Item(19, 16).Value = 35  # P19 n
$ws.Cells.Item(19, 17).Value = 7  # Q19 p_typicality
$ws.Cells.Item(19, 18).Value = 7  # R19 p_conceptual
$ws.Cells.Item(19, 19).Value = 7  # S19 p_perceptual

# row 20 (trial_block 19)
$ws.Cells.Item(20, 6).Value = 207  # F20 trial_total
$ws.Cells.Item(20, 9).Value = $null  # I20 cond_cat
$ws.Cells.Item(20, 10).Value = 'new'  # J20 cond_mem
$ws.Cells.Item(20, 11).Value = 'f'  # K20 correct_answer
$ws.Cells.Item(20, 12).Value = 'stimuli/img_bklr1.png'  # L20 stimulus
$ws.Cells.Item(20, 13).Value = 86.54761904761905  # M20 conceptual
$ws.Cells.Item(20, 14).Value = 67.73809523809524  # N20 perceptual
$ws.Cells.Item(20, 15).Value = 77.14285714285714  # O20 typicality
$ws.Cells.Item(20, 16).Value = 42  # P20 n
$ws.Cells.Item(20, 17).Value = 9  # Q20 p_typicality
$ws.Cells.Item(20, 18).Value = 9  # R20 p_conceptual
$ws.Cells.Item(20, 19).Value = 9  # S20 p_perceptual

# row 21 (trial_block 20)
$ws.Cells.Item(21, 6).Value = 208  # F21 trial_total
$ws.Cells.Item(21, 12).Value = 'stimuli/img_2pk6v.png'  # L21 stimulus
$ws.Cells.Item(21, 13).Value = 85.08108108108108  # M21 conceptual
$ws.Cells.Item(21, 14).Value = 66.16216216216216  # N21 perceptual
$ws.Cells.Item(21, 15).Value = 75.62162162162161  # O21 typicality
$ws.Cells.Item(21, 17).Value = 9  # Q21 p_typicality
$ws.Cells.Item(21, 18).Value = 9  # R21 p_conceptual
$ws.Cells.Item(21, 19).Value = 9  # S21 p_perceptual

# row 22 (trial_block 21)
$ws.Cells.Item(22, 6).Value = 209  # F22 trial_total
$ws.Cells.Item(22, 9).Value = 'target'  # I22 cond_cat
$ws.Cells.Item(22, 10).Value = 'old'  # J22 cond_mem
$ws.Cells.Item(22, 11).Value = 'j'  # K22 correct_answer
$ws.Cells.Item(22, 12).Value = 'stimuli/img_h0hbk.png'  # L22 stimulus
$ws.Cells.Item(22, 13).Value = 86.80952380952381  # M22 conceptual
$ws.Cells.Item(22, 14).Value = 69.19047619047619  # N22 perceptual
$ws.Cells.Item(22, 15).Value = 78  # O22 typicality
$ws.Cells.Item(22, 16).Value = 42  # P22 n
$ws.Cells.Item(22, 17).Value = 9  # Q22 p_typicality
$ws.Cells.Item(22, 18).Value = 9  # R22 p_conceptual
$ws.Cells.Item(22, 19).Value = 9  # S22 p_perceptual

# row 23 (trial_block 22)
$ws.Cells.Item(23, 6).Value = 210  # F23 trial_total
$ws.Cells.Item(23, 12).Value = 'stimuli/img_x0u5z.png'  # L23 stimulus
$ws.Cells.Item(23, 13).Value = 92  # M23 conceptual
$ws.Cells.Item(23, 14).Value = 78.16216216216216  # N23 perceptual
$ws.Cells.Item(23, 15).Value = 85.08108108108108  # O23 typicality

# row 24 (trial_block 23)
$ws.Cells.Item(24, 6).Value = 211  # F24 trial_total
$ws.Cells.Item(24, 9).Value = $null  # I24 cond_cat
$ws.Cells.Item(24, 10).Value = 'new'  # J24 cond_mem
$ws.Cells.Item(24, 11).Value = 'f'  # K24 correct_answer
$ws.Cells.Item(24, 12).Value = 'stimuli/img_5yhyk.png'  # L24 stimulus
$ws.Cells.Item(24, 13).Value = 46.375  # M24 conceptual
$ws.Cells.Item(24, 14).Value = 31.325  # N24 perceptual
$ws.Cells.Item(24, 15).Value = 38.85  # O24 typicality
$ws.Cells.Item(24, 16).Value = 40  # P24 n
$ws.Cells.Item(24, 17).Value = 2  # Q24 p_typicality
$ws.Cells.Item(24, 18).Value = 2  # R24 p_conceptual
$ws.Cells.Item(24, 19).Value = 2  # S24 p_perceptual

# row 25 (trial_block 24)
$ws.Cells.Item(25, 6).Value = 212  # F25 trial_total
$ws.Cells.Item(25, 9).Value = $null  # I25 cond_cat
$ws.Cells.Item(25, 10).Value = 'new'  # J25 cond_mem
$ws.Cells.Item(25, 11).Value = 'f'  # K25 correct_answer
$ws.Cells.Item(25, 12).Value = 'stimuli/img_i7vab.png'  # L25 stimulus
$ws.Cells.Item(25, 13).Value = 86.40000000000001  # M25 conceptual
$ws.Cells.Item(25, 14).Value = 67.8  # N25 perceptual
$ws.Cells.Item(25, 15).Value = 77.09999999999999  # O25 typicality
$ws.Cells.Item(25, 16).Value = 35  # P25 n
$ws.Cells.Item(25, 17).Value = 9  # Q25 p_typicality
$ws.Cells.Item(25, 18).Value = 9  # R25 p_conceptual
$ws.Cells.Item(25, 19).Value = 9  # S25 p_perceptual

# row 26 (trial_block 25)
$ws.Cells.Item(26, 6).Value = 213  # F26 trial_total

# row 27 (trial_block 26)
$ws.Cells.Item(27, 6).Value = 214  # F27 trial_total
$ws.Cells.Item(27, 12).Value = 'stimuli/img_t2ioc.png'  # L27 stimulus
$ws.Cells.Item(27, 13).Value = 88.18918918918919  # M27 conceptual
$ws.Cells.Item(27, 14).Value = 74.05405405405405  # N27 perceptual
$ws.Cells.Item(27, 15).Value = 81.12162162162161  # O27 typicality
$ws.Cells.Item(27, 16).Value = 37  # P27 n
$ws.Cells.Item(27, 17).Value = 10  # Q27 p_typicality
$ws.Cells.Item(27, 18).Value = 10  # R27 p_conceptual
$ws.Cells.Item(27, 19).Value = 10  # S27 p_perceptual

# row 28 (trial_block 27)
$ws.Cells.Item(28, 6).Value = 215  # F28 trial_total
$ws.Cells.Item(28, 12).Value = 'stimuli/img_ybbmx.png'  # L28 stimulus
$ws.Cells.Item(28, 13).Value = 55.24324324324324  # M28 conceptual
$ws.Cells.Item(28, 14).Value = 36.75675675675676  # N28 perceptual
$ws.Cells.Item(28, 15).Value = 46  # O28 typicality
$ws.Cells.Item(28, 16).Value = 37  # P28 n
$ws.Cells.Item(28, 17).Value = 3  # Q28 p_typicality
$ws.Cells.Item(28, 18).Value = 3  # R28 p_conceptual
$ws.Cells.Item(28, 19).Value = 3  # S28 p_perceptual

# row 29 (trial_block 28)
$ws.Cells.Item(29, 6).Value = 216  # F29 trial_total
